$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: height change only (becomes an explicit custom height)
$ws.Rows.Item(12).RowHeight = 43.6

# Row 13: restyle existing cells (vertical=center, wrapText=true) and add F13
$ws.Range("A13").VerticalAlignment = -4108
$ws.Range("B13").VerticalAlignment = -4108
$ws.Range("D13").VerticalAlignment = -4108

$ws.Range("C13").VerticalAlignment = -4108
$ws.Range("C13").WrapText = $true
$ws.Range("E13").VerticalAlignment = -4108
$ws.Range("E13").WrapText = $true

$ws.Range("F13").VerticalAlignment = -4108
$ws.Range("F13").WrapText = $true

$ws.Rows.Item(13).RowHeight = 35.55

# Row 14: new data row
$ws.Range("A14").Value = 11
$ws.Range("B14").Value = "Escopo"
$ws.Range("C14").Value = "Não estava explicito no escopo que a produção do chip e leitor seriam feito de forma terceirizada. "
$ws.Range("D14").Value = "Escopo"
$ws.Range("E14").Value = "Foi incluído uma restrição da declaração de escopo."

$ws.Range("A14").VerticalAlignment = -4108
$ws.Range("B14").VerticalAlignment = -4108
$ws.Range("D14").VerticalAlignment = -4108

$ws.Range("C14").VerticalAlignment = -4108
$ws.Range("C14").WrapText = $true
$ws.Range("E14").VerticalAlignment = -4108
$ws.Range("E14").WrapText = $true

$ws.Range("F14").VerticalAlignment = -4108
$ws.Range("F14").WrapText = $true

$ws.Rows.Item(14).RowHeight = 59.7

# Selection / view
[void]$ws.Range("F11").Select()

Write-Output "done"
